$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($ws, $addr, $val)
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-CellText $ws 'D2' '26.893.62'
Set-CellText $ws 'E2' '  +0.13%  '
Set-CellText $ws 'D3' '1.545.20'
Set-CellText $ws 'E3' '  -1.14%  '
Set-CellText $ws 'E4' '  +0.22%  '
Set-CellText $ws 'D5' '205.85'
Set-CellText $ws 'E5' '  -0.05%  '
Set-CellText $ws 'E6' '  -0.48%  '
Set-CellText $ws 'E8' '  -0.13%  '
Set-CellText $ws 'D9' '21.24'
Set-CellText $ws 'E9' '  -2.22%  '
Set-CellText $ws 'D10' '0.0582'
Set-CellText $ws 'D11' '0.0857'
Set-CellText $ws 'E11' '  -0.83%  '
Set-CellText $ws 'D12' '1.764.53'
Set-CellText $ws 'E12' '  -1.14%  '
Set-CellText $ws 'D13' '1.547.51'
Set-CellText $ws 'E13' '  -0.99%  '
Set-CellText $ws 'D14' '3.70'
Set-CellText $ws 'E14' '  -0.81%  '
Set-CellText $ws 'E15' '  -0.76%  '
Set-CellText $ws 'D16' '26.877.05'
Set-CellText $ws 'E16' '  +0.05%  '
Set-CellText $ws 'D17' '61.39'
Set-CellText $ws 'E17' '  +0.20%  '
Set-CellText $ws 'D18' '213.44'
Set-CellText $ws 'E18' '  -0.75%  '
Set-CellText $ws 'D19' '0.0₃0681'
Set-CellText $ws 'E19' '  +0.19%  '
Set-CellText $ws 'D20' '7.16'
Set-CellText $ws 'E20' '  -2.49%  '
Set-CellText $ws 'E21' '  +0.28%  '
Set-CellText $ws 'E22' '  -2.61%  '
Set-CellText $ws 'D23' '9.18'
Set-CellText $ws 'E23' '  +0.21%  '
Set-CellText $ws 'D24' '1.94'
Set-CellText $ws 'E24' '  -3.24%  '
Set-CellText $ws 'D25' '152.37'
Set-CellText $ws 'E25' '  -0.99%  '
Set-CellText $ws 'D26' '6.60'
Set-CellText $ws 'E26' '  -1.87%  '
Set-CellText $ws 'D27' '14.80'
Set-CellText $ws 'E27' '  -0.75%  '
Set-CellText $ws 'E28' '  +0.33%  '
Set-CellText $ws 'E29' '  -0.09%  '
Set-CellText $ws 'D30' '0.0458'
Set-CellText $ws 'E30' '  -1.55%  '
Set-CellText $ws 'E31' '  -1.02%  '
Set-CellText $ws 'E32' '  +1.46%  '
Set-CellText $ws 'D33' '1.355.88'
Set-CellText $ws 'E33' '  -3.43%  '
Set-CellText $ws 'D34' '2.93'
Set-CellText $ws 'E34' '  +0.47%  '
Set-CellText $ws 'E35' '  +0.25%  '
Set-CellText $ws 'D36' '0.963'
Set-CellText $ws 'E36' '  +5.11%  '
Set-CellText $ws 'D37' '2.28'
Set-CellText $ws 'E37' '  +0.34%  '
Set-CellText $ws 'E38' '  -0.32%  '
Set-CellText $ws 'D39' '0.519'
Set-CellText $ws 'E39' '  -1.75%  '
Set-CellText $ws 'D40' '0.803'
Set-CellText $ws 'E40' '  -0.85%  '
Set-CellText $ws 'E41' '  +0.31%  '
Set-CellText $ws 'D42' '5.55'
Set-CellText $ws 'E42' '  +1.98%  '
Set-CellText $ws 'D43' '0.988'
Set-CellText $ws 'E43' '  -0.91%  '
Set-CellText $ws 'E44' '  +1.95%  '
Set-CellText $ws 'D45' '63.26'
Set-CellText $ws 'E45' '  -0.02%  '
Set-CellText $ws 'D46' '1.72'
Set-CellText $ws 'E46' '  -2.14%  '
Set-CellText $ws 'D47' '1.678.56'
Set-CellText $ws 'E47' '  -1.22%  '
Set-CellText $ws 'B48' 'mCoin'
Set-CellText $ws 'C48' 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
Set-CellText $ws 'D48' '2.21'
Set-CellText $ws 'E48' '  -5.62%  '
Set-CellText $ws 'B49' 'Quant'
Set-CellText $ws 'C49' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-CellText $ws 'D49' '85.74'
Set-CellText $ws 'E49' '  -0.63%  '
Set-CellText $ws 'B50' 'Cronos'
Set-CellText $ws 'C50' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-CellText $ws 'D50' '0.0510'
Set-CellText $ws 'E50' '  +0.92%  '
Set-CellText $ws 'B51' 'BabyDogeCoin'
Set-CellText $ws 'C51' 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-CellText $ws 'D51' '0.0₇0975'
Set-CellText $ws 'E51' '  -0.37%  '
